$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update totals affected by the new period being added ---
# VALOR MORA total (was 604266, +52000 for the new 2509 period)
$ws.Range("E11").Value = 656266
# Cant. Periodos (was 9, +1 for the new 2509 period)
$ws.Range("F13").Value = 10

# --- Insert a new data row for period "2509" right after the last existing
#     data row (25), pushing the signature block rows down by one ---
$ws.Rows("26").Insert()

# The newly inserted row 26 starts out blank/unformatted. Give it the same
# "closing" look (bottom border, etc.) that the previous last row (25) had,
# by copying that row's formatting down into it.
$ws.Range("B25:J25").Copy()
$ws.Range("B26:J26").PasteSpecial(-4122)

# Row 25 is no longer the last row of the table, so restore its formatting
# to match the regular (non-closing) rows above it (16-24).
$ws.Range("B24:J24").Copy()
$ws.Range("B25:J25").PasteSpecial(-4122)

# Fill in the data for the new period row
$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1047386377"
$ws.Range("D26").Value = "EDGAR JOSE ROA AMADOR"
$ws.Range("E26").Value = "2509"
$ws.Range("F26").Value = 52000
$ws.Range("G26").Value = 1300000

# Center the "Periodo Mora" column for the two bottom rows
$ws.Range("E25").HorizontalAlignment = -4108
$ws.Range("E26").HorizontalAlignment = -4108

$excel.CutCopyMode = 0
